$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-11: 45183 -> 45184
$ws.Range("C2:C11").Value = 45184

# Row 2 hyperlink formulas: add the link-text second argument.
# S2 keeps the (buggy) quoting exactly as committed upstream - the
# closing quote of the URL string was placed after the new text
# instead of before it.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/artfynd/A 33491-2023.xlsx, "A 33491-2023"")'

$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/kartor/A 33491-2023.png", "A 33491-2023")'

$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/klagomål/A 33491-2023.docx", "A 33491-2023")'

$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/klagomålsmail/A 33491-2023.docx", "A 33491-2023")'

$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/tillsyn/A 33491-2023.docx", "A 33491-2023")'

# Y2 used to be a literal inline-string (not a real formula, with a
# Swedish ";" argument separator). Replace it with a proper formula
# using the standard "," separator, matching the other link columns.
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/tillsynsmail/A 33491-2023.docx", "A 33491-2023")'
